# Regenerate merged AHB files
# - rename the "_old" / "_new" header-suffix columns to "_FV2404" / "_FV2410"
# - freeze the header row
# - wrap the data range in an Excel Table (Table1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename header row labels (row 1): "_old" -> "_FV2404", "_new" -> "_FV2410"
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

# Columns A-J -> "<name>_FV2404"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2404"
}

# Column K stays "diff" (unchanged)

# Columns L-U -> "<name>_FV2410"
for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = 12 + $i
    $ws.Cells.Item(1, $col).Value = "$($baseNames[$i])_FV2410"
}

# ---------------------------------------------------------------------------
# 2) Freeze the header row (pane split after row 1)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3) Turn the data range into an Excel Table ("Table1")
# ---------------------------------------------------------------------------
$tableRange = $ws.Range("A1:U73")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"
